# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - mirror the look of the existing header cells (A1:AC1):
# bold font, thin border all around, centered horizontally, top-aligned vertically.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop

# Season record is the same for every player row (2-52)
$lastRow = 52
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 71   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 91   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
